$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = "Riccardo Barbiero"
$ws.Range("B66").Value = "Elia Battisti | U.SGUARNA"
$ws.Range("C66").Value = "Riccardo Barbiero | Rita Levi’s"
$ws.Range("D66").Value = "Leonardo Viola | SHARK ATTACK"
$ws.Range("E66").Value = "Marco Sala | IMONTAGNA"
$ws.Range("F66").Value = "Moris Benedetti | Gli Introvabili"
